$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Experimental value: set to the literal text "false" -----------------
# Assigning the bare string "false" gets auto-coerced to the Boolean FALSE
# by the engine, so we stage it with a leading space (kept as text), then
# use a helper cell with TRIM() copied back as a value to normalize it to
# the exact text "false" without leaving a formula or changing styles.
$ws.Cells.Item(7, 2).Value = " false"
$helper = $ws.Cells.Item(500, 2)
$helper.Formula = "=TRIM(B7)"
$helper.Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4163)  # xlPasteValues
$helper.Value = ""

# --- Date value update -----------------------------------------------------
$ws.Cells.Item(8, 2).Value = "2025-11-30T13:08:37+00:00"

# --- Description value ------------------------------------------------------
$ws.Cells.Item(17, 2).Value = "Exercise protocols used for VO2max testing"

$excel.CutCopyMode = 0
